$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must remain text even though it looks like a plain
# decimal number (Excel would otherwise auto-convert it to a Number). We force the
# "Text" number format just long enough to type the value in as a string, then strip
# the formatting back off so the cell keeps its original (default) style - only the
# underlying stored type (text) changes, matching the source inline-string cells.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

$ws.Range('D2').Value = '39.690.48'
$ws.Range('E2').Value = '  +0.44%  '
$ws.Range('D3').Value = '2.173.80'
$ws.Range('E4').Value = '  +0.12%  '
Set-TextValue $ws.Range('D5') '226.58'
$ws.Range('E5').Value = '  -1.19%  '
$ws.Range('E6').Value = '  +0.64%  '
Set-TextValue $ws.Range('D7') '63.01'
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('E9').Value = '  -0.70%  '
Set-TextValue $ws.Range('D10') '0.0850'
$ws.Range('E10').Value = '  -1.21%  '
$ws.Range('E11').Value = '  +0.38%  '
$ws.Range('E12').Value = '  -1.48%  '
$ws.Range('D13').Value = '2.495.72'
$ws.Range('E13').Value = '  +0.70%  '
$ws.Range('E14').Value = '  -2.16%  '
$ws.Range('E15').Value = '  -1.44%  '
$ws.Range('E16').Value = '  -1.45%  '
$ws.Range('D17').Value = '2.179.46'
$ws.Range('E17').Value = '  -0.50%  '
$ws.Range('D18').Value = '39.654.71'
$ws.Range('E18').Value = '  +0.43%  '
$ws.Range('D19').Value = '0.0₃0914'
$ws.Range('E19').Value = '  +6.91%  '
$ws.Range('E20').Value = '  -0.86%  '
Set-TextValue $ws.Range('D21') '6.01'
$ws.Range('E21').Value = '  -2.36%  '
Set-TextValue $ws.Range('D22') '229.57'
$ws.Range('E22').Value = '  +0.45%  '
$ws.Range('E24').Value = '  -0.89%  '
Set-TextValue $ws.Range('D25') '2.34'
$ws.Range('E25').Value = '  -1.00%  '
Set-TextValue $ws.Range('D26') '9.55'
$ws.Range('E26').Value = '  -2.35%  '
Set-TextValue $ws.Range('D27') '170.88'
$ws.Range('E27').Value = '  -1.25%  '
Set-TextValue $ws.Range('D28') '0.140'
$ws.Range('E28').Value = '  +0.10%  '
$ws.Range('E29').Value = '  +1.13%  '
Set-TextValue $ws.Range('D30') '19.81'
$ws.Range('E30').Value = '  +0.77%  '
Set-TextValue $ws.Range('D31') '2.68'
$ws.Range('E31').Value = '  +4.08%  '
$ws.Range('E32').Value = '  +0.15%  '
$ws.Range('E33').Value = '  -2.55%  '
Set-TextValue $ws.Range('D34') '4.69'
$ws.Range('E34').Value = '  -2.75%  '
$ws.Range('E35').Value = '  -2.46%  '
$ws.Range('E36').Value = '  -0.85%  '
$ws.Range('E37').Value = '  +8.30%  '
$ws.Range('E38').Value = '  -1.90%  '
$ws.Range('E39').Value = '  +0.18%  '
Set-TextValue $ws.Range('D40') '4.91'
$ws.Range('E40').Value = '  +12.65%  '
$ws.Range('E41').Value = '  -0.87%  '
Set-TextValue $ws.Range('D42') '102.54'
$ws.Range('E42').Value = '  -0.82%  '
Set-TextValue $ws.Range('D43') '17.70'
$ws.Range('E43').Value = '  -2.65%  '
$ws.Range('D44').Value = '1.512.14'
$ws.Range('E44').Value = '  -1.24%  '
Set-TextValue $ws.Range('D45') '1.22'
$ws.Range('E45').Value = '  +2.38%  '
Set-TextValue $ws.Range('D46') '7.90'
$ws.Range('E46').Value = '  +1.74%  '
$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws.Range('D47') '0.0920'
$ws.Range('E47').Value = '  -0.96%  '
$ws.Range('B48').Value = 'HuobiToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue $ws.Range('D48') '2.79'
$ws.Range('E48').Value = '  -0.54%  '
$ws.Range('E49').Value = '  -1.55%  '
Set-TextValue $ws.Range('D50') '0.000195'
$ws.Range('E50').Value = '  +33.06%  '
Set-TextValue $ws.Range('D51') '49.38'
$ws.Range('E51').Value = '  +6.04%  '
